$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 51.59157666666666
$ws.Range("H2").Value = 154.77473
$ws.Range("I2").Value = 0.2641250550177587
$ws.Range("J2").Value = 0.2641250550177588
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.3116673333333334
$ws.Range("N2").Value = 0.935002
$ws.Range("O2").Value = 0.0414413620607491
$ws.Range("P2").Value = 0.0414413620607491
$ws.Range("Q2").Value = 16.07940912216222
$ws.Range("R2").Value = 144.71468209946
$ws.Range("S2").Value = 0.01094570203430621
$ws.Range("T2").Value = 0.01094570203430622

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 51.59157666666666
$ws.Range("H3").Value = 154.77473
$ws.Range("I3").Value = 0.2641250550177587
$ws.Range("J3").Value = 0.2641250550177588
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.794584666666667
$ws.Range("N3").Value = 11.383754
$ws.Range("O3").Value = 0.5045532214096876
$ws.Range("P3").Value = 0.5045532214096876
$ws.Range("Q3").Value = 195.7686057484911
$ws.Range("R3").Value = 1761.91745173642
$ws.Range("S3").Value = 0.1332651473642212
$ws.Range("T3").Value = 0.1332651473642212

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 51.59157666666666
$ws.Range("H4").Value = 154.77473
$ws.Range("I4").Value = 0.2641250550177587
$ws.Range("J4").Value = 0.2641250550177588
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.414430666666667
$ws.Range("N4").Value = 10.243292
$ws.Range("O4").Value = 0.4540054165295632
$ws.Range("P4").Value = 0.4540054165295633
$ws.Range("Q4").Value = 176.1558615123511
$ws.Range("R4").Value = 1585.40275361116
$ws.Range("S4").Value = 0.1199142056192314
$ws.Range("T4").Value = 0.1199142056192314

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.32115333333334
$ws.Range("H5").Value = 57.96346000000001
$ws.Range("I5").Value = 0.09891538535728452
$ws.Range("J5").Value = 0.09891538535728453
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.3116673333333334
$ws.Range("N5").Value = 0.935002
$ws.Range("O5").Value = 0.0414413620607491
$ws.Range("P5").Value = 0.0414413620607491
$ws.Range("Q5").Value = 6.021772336324446
$ws.Range("R5").Value = 54.19595102692001
$ws.Range("S5").Value = 0.004099188297969747
$ws.Range("T5").Value = 0.004099188297969748

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.32115333333334
$ws.Range("H6").Value = 57.96346000000001
$ws.Range("I6").Value = 0.09891538535728452
$ws.Range("J6").Value = 0.09891538535728453
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.794584666666667
$ws.Range("N6").Value = 11.383754
$ws.Range("O6").Value = 0.5045532214096876
$ws.Range("P6").Value = 0.5045532214096876
$ws.Range("Q6").Value = 73.31575218098224
$ws.Range("R6").Value = 659.8417696288401
$ws.Range("S6").Value = 0.04990807632899855
$ws.Range("T6").Value = 0.04990807632899855

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.32115333333334
$ws.Range("H7").Value = 57.96346000000001
$ws.Range("I7").Value = 0.09891538535728452
$ws.Range("J7").Value = 0.09891538535728453
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.414430666666667
$ws.Range("N7").Value = 10.243292
$ws.Range("O7").Value = 0.4540054165295632
$ws.Range("P7").Value = 0.4540054165295633
$ws.Range("Q7").Value = 65.97073845670224
$ws.Range("R7").Value = 593.7366461103202
$ws.Range("S7").Value = 0.04490812073031621
$ws.Range("T7").Value = 0.04490812073031623

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 112.3724673333333
$ws.Range("H8").Value = 337.117402
$ws.Range("I8").Value = 0.5752951554216499
$ws.Range("J8").Value = 0.57529515542165
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.3116673333333334
$ws.Range("N8").Value = 0.935002
$ws.Range("O8").Value = 0.0414413620607491
$ws.Range("P8").Value = 0.0414413620607491
$ws.Range("Q8").Value = 35.02282723386711
$ws.Range("R8").Value = 315.205445104804
$ws.Range("S8").Value = 0.02384101482762352
$ws.Range("T8").Value = 0.02384101482762352

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 112.3724673333333
$ws.Range("H9").Value = 337.117402
$ws.Range("I9").Value = 0.5752951554216499
$ws.Range("J9").Value = 0.57529515542165
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.794584666666667
$ws.Range("N9").Value = 11.383754
$ws.Range("O9").Value = 0.5045532214096876
$ws.Range("P9").Value = 0.5045532214096876
$ws.Range("Q9").Value = 426.4068414985675
$ws.Range("R9").Value = 3837.661573487107
$ws.Range("S9").Value = 0.2902670239293804
$ws.Range("T9").Value = 0.2902670239293804

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 112.3724673333333
$ws.Range("H10").Value = 337.117402
$ws.Range("I10").Value = 0.5752951554216499
$ws.Range("J10").Value = 0.57529515542165
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.414430666666667
$ws.Range("N10").Value = 10.243292
$ws.Range("O10").Value = 0.4540054165295632
$ws.Range("P10").Value = 0.4540054165295633
$ws.Range("Q10").Value = 383.6879985519315
$ws.Range("R10").Value = 3453.191986967384
$ws.Range("S10").Value = 0.2611871166646459
$ws.Range("T10").Value = 0.2611871166646461

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 12.044915
$ws.Range("H11").Value = 36.134745
$ws.Range("I11").Value = 0.06166440420330686
$ws.Range("J11").Value = 0.06166440420330688
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.3116673333333334
$ws.Range("N11").Value = 0.935002
$ws.Range("O11").Value = 0.0414413620607491
$ws.Range("P11").Value = 0.0414413620607491
$ws.Range("Q11").Value = 3.754006538276667
$ws.Range("R11").Value = 33.78605884449
$ws.Range("S11").Value = 0.002555456900849618
$ws.Range("T11").Value = 0.002555456900849619

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 12.044915
$ws.Range("H12").Value = 36.134745
$ws.Range("I12").Value = 0.06166440420330686
$ws.Range("J12").Value = 0.06166440420330688
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.794584666666667
$ws.Range("N12").Value = 11.383754
$ws.Range("O12").Value = 0.5045532214096876
$ws.Range("P12").Value = 0.5045532214096876
$ws.Range("Q12").Value = 45.70544977030334
$ws.Range("R12").Value = 411.34904793273
$ws.Range("S12").Value = 0.03111297378708756
$ws.Range("T12").Value = 0.03111297378708757

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 12.044915
$ws.Range("H13").Value = 36.134745
$ws.Range("I13").Value = 0.06166440420330686
$ws.Range("J13").Value = 0.06166440420330688
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.414430666666667
$ws.Range("N13").Value = 10.243292
$ws.Range("O13").Value = 0.4540054165295632
$ws.Range("P13").Value = 0.4540054165295633
$ws.Range("Q13").Value = 176.1558615123511
$ws.Range("R13").Value = 1585.40275361116
$ws.Range("S13").Value = 0.1199142056192314
$ws.Range("T13").Value = 0.1199142056192314
